# Applies:
#  1) The table on slide 16 switches from table-style {F0008912-6605-40DE-BC96-891BA60F547A}
#     to {1A9674B3-44E6-4BF8-B7E7-91EB20A745E2}.
#  2) The presentation's theme colour scheme (the "Integral" design applied to the
#     slide master) is changed back to the stock "Office" colour palette.

function ConvertTo-VbaRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 ------------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{1A9674B3-44E6-4BF8-B7E7-91EB20A745E2}")
    }
}

# --- 2) Theme colours: Integral -> Office --------------------------------------
# ThemeColorScheme.Item(index) order: dk1, lt1, dk2, lt2,
# accent1..accent6, hlink, folHlink
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = ConvertTo-VbaRGB $officeColors[$i - 1]
}
